# Commit: "change metadata sheet to isa template"
# The Swate metadata worksheet is renamed from "SwateTemplateMetadata"
# to "isa_template".

$wb = $excel.ActiveWorkbook

$metadataSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metadataSheet.Name = "isa_template"
